# Update "Estado de Cuenta" worksheet: refresh worker/debt data (part 1 of new batch)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Header / summary cells
# ---------------------------------------------------------------------------
$ws.Range("D2").Value2  = "ESTADO DE CUENTA"
$ws.Range("B7").Value2  = "RAZON SOCIAL:"
$ws.Range("B11").Value2 = "VALOR MORA"
$ws.Range("E11").Value2 = 504868

$ws.Range("B13").Value2 = "Cant. Trabajadores"
$ws.Range("C13").Value2 = 14
$ws.Range("E13").Value2 = "Cant. Periodos"
$ws.Range("F13").Value2 = 1

$ws.Range("H15").Value2 = "Novedad de Ingreso"
$ws.Range("I15").Value2 = "Novedad de Retiro"
$ws.Range("J15").Value2 = "Observaciones"

# ---------------------------------------------------------------------------
# 2) Make room for the extra worker rows.
#    Before: 6 debt rows (16-21), last one (21) carries the "closing" bottom
#    border. After: 14 debt rows (16-29), so insert 8 blank rows right below
#    the current last row (21) - this pushes the old signature block
#    (rows 26-27) down to rows 34-35 automatically.
# ---------------------------------------------------------------------------
$ws.Rows("22:29").Insert()

# Copy the format of a normal (non-closing) row into the 7 brand-new interior
# rows (22-28).
$ws.Range("B20:J20").Copy()
$ws.Range("B22:J28").PasteSpecial(-4122)  # xlPasteFormats

# The old row 21 still carries the "closing" (special) border - move that
# format down onto the new last row (29).
$ws.Range("B21:J21").Copy()
$ws.Range("B29:J29").PasteSpecial(-4122)  # xlPasteFormats

# Row 21 itself becomes a normal interior row now.
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Worker / debt-period rows (16-29): Tipo Doc, N Doc, Nombre, Periodo,
#    Valor Mora, Salario Basico. H/I/J (Novedad Ingreso/Retiro/Observaciones)
#    stay blank, like before.
# ---------------------------------------------------------------------------
$rows = @(
  @{r=16; doc="1104011715"; nombre="RAFAEL ENRIQUE DAVILA CASTILLO"},
  @{r=17; doc="19873532";   nombre="JORGE LUIS MOZO CAEZ"},
  @{r=18; doc="1099965448"; nombre="JORGE LUIS MOZO GAVIRIA"},
  @{r=19; doc="18762742";   nombre="SERGIO JOSE RIVERA CARPIO"},
  @{r=20; doc="18761643";   nombre="JOSE DANIEL PEREZ RAMIREZ"},
  @{r=21; doc="18762463";   nombre="CARLOS ANDRES ARRIETA ROMERO"},
  @{r=22; doc="92189411";   nombre="ALDEMAR JOSE GAMBOA CAMPO"},
  @{r=23; doc="1100625745"; nombre="ANDERSON DAVID GARCIA ROMERO"},
  @{r=24; doc="1193113822"; nombre="DIEGO ARMANDO BENITEZ VILLEGAS"},
  @{r=25; doc="92188565";   nombre="NIBALDO RAUL SARMIENTO BARRIOS"},
  @{r=26; doc="9218901";    nombre="JAIRO ALBERTO MARTINEZ SANTOS"},
  @{r=27; doc="92190046";   nombre="ALFREDO JOSE ROMERO MERCADO"},
  @{r=28; doc="92191201";   nombre="WILFRAN ANTONIO SOLORZANO BOHORQUEZ"},
  @{r=29; doc="92188901";   nombre="JAIRO ALBERTO MARTINEZ SANTIZ"}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Range("B$r").Value2 = "CC"
  $ws.Range("C$r").Value2 = $row.doc
  $ws.Range("D$r").Value2 = $row.nombre
  $ws.Range("E$r").Value2 = "2508"
  $ws.Range("F$r").Value2 = 36062
  $ws.Range("G$r").Value2 = 1423500
}

# ---------------------------------------------------------------------------
# 4) Widen column D so the longer names fit ("bestFit" already recorded by
#    Excel at 41.81640625 once data is in place).
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 41.81640625
